$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.205.52"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.070.35"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'231.54"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'58.11"
$ws.Range("E7").Value = "  +5.32%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.388"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "2.364.81"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "'14.63"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "'20.76"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.753"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "2.058.81"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "38.050.72"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "'6.28"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'70.13"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "'225.64"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "'2.28"
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("D26").Value = "'9.32"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'166.13"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'0.135"
$ws.Range("E28").Value = "  +6.29%  "
$ws.Range("D29").Value = "'19.17"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("D30").Value = "'1.37"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'0.119"
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("D32").Value = "'4.57"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'4.60"
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("D35").Value = "'1.99"
$ws.Range("E35").Value = "  +7.64%  "
$ws.Range("D36").Value = "'2.39"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").Value = "'6.00"
$ws.Range("E37").Value = "  +10.82%  "
$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'98.84"
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("D41").Value = "'0.0220"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "1.486.09"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'0.0955"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'16.86"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.85"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").Value = "'1.13"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").Value = "'4.07"
$ws.Range("E47").Value = "  +14.26%  "
$ws.Range("D48").Value = "'1.03"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.96"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'7.12"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").Value = "2.247.95"
$ws.Range("E51").Value = "  +1.65%  "
